$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,12
$row2[0,0] = 0.1092194802715341
$row2[0,1] = 0.1199478036574728
$row2[0,2] = 0.1178313381864697
$row2[0,3] = 2.323063226037007
$row2[0,4] = 0.002547691809639624
$row2[0,5] = 0
$row2[0,6] = 2.175360659233633
$row2[0,7] = 0
$row2[0,8] = 2.965446372802319
$row2[0,9] = 0.1921021294112393
$row2[0,10] = 0
$row2[0,11] = 1.683160507300052
$ws.Range("C2:N2").Value = $row2

$row3 = New-Object 'object[,]' 1,12
$row3[0,0] = 0.1068530955642188
$row3[0,1] = 0.1214616737089962
$row3[0,2] = 0.1151583143973767
$row3[0,3] = 2.299663445876291
$row3[0,4] = 0.002553447920325986
$row3[0,5] = 0
$row3[0,6] = 2.158183152063543
$row3[0,7] = 0
$row3[0,8] = 2.785727557391112
$row3[0,9] = 0.1876275752744618
$row3[0,10] = 0
$row3[0,11] = 1.706726430940291
$ws.Range("C3:N3").Value = $row3

$row4 = New-Object 'object[,]' 1,12
$row4[0,0] = 0.1054477437219674
$row4[0,1] = 0.1224413576291159
$row4[0,2] = 0.1135841604944439
$row4[0,3] = 2.286878310438993
$row4[0,4] = 0.002557165638417856
$row4[0,5] = 0
$row4[0,6] = 2.149060192161443
$row4[0,7] = 0
$row4[0,8] = 2.676977010506391
$row4[0,9] = 0.1850014338679529
$row4[0,10] = 0
$row4[0,11] = 1.721881643665386
$ws.Range("C4:N4").Value = $row4

$row5 = New-Object 'object[,]' 1,12
$row5[0,0] = 0.1048870325352169
$row5[0,1] = 0.1228531700957625
$row5[0,2] = 0.1129594901215789
$row5[0,3] = 2.282064090891097
$row5[0,4] = 0.002558726932609514
$row5[0,5] = 0
$row5[0,6] = 2.145698892959103
$row5[0,7] = 0
$row5[0,8] = 2.633059511302804
$row5[0,9] = 0.183961587352627
$row5[0,10] = 0
$row5[0,11] = 1.728229614403321
$ws.Range("C5:N5").Value = $row5

$row6 = New-Object 'object[,]' 1,12
$row6[0,0] = 0.1047946507385689
$row6[0,1] = 0.1229223103715498
$row6[0,2] = 0.1128567776189087
$row6[0,3] = 2.281288539843771
$row6[0,4] = 0.002558988985378389
$row6[0,5] = 0
$row6[0,6] = 2.145162226443318
$row6[0,7] = 0
$row6[0,8] = 2.625791085590379
$row6[0,9] = 0.1837907486356798
$row6[0,10] = 0
$row6[0,11] = 1.729294075460838
$ws.Range("C6:N6").Value = $row6

$row7 = New-Object 'object[,]' 1,12
$row7[0,0] = 0.1054401332526851
$row7[0,1] = 0.122446860567651
$row7[0,2] = 0.113575667986467
$row7[0,3] = 2.28681178428009
$row7[0,4] = 0.002557186507175248
$row7[0,5] = 0
$row7[0,6] = 2.149013419822651
$row7[0,7] = 0
$row7[0,8] = 2.676383110528832
$row7[0,9] = 0.1849872875504843
$row7[0,10] = 0
$row7[0,11] = 1.721966558230655
$ws.Range("C7:N7").Value = $row7

$row8 = New-Object 'object[,]' 1,12
$row8[0,0] = 0.1083936707052686
$row8[0,1] = 0.1204593568168733
$row8[0,2] = 0.1168957160322712
$row8[0,3] = 2.314665273467725
$row8[0,4] = 0.00254963854259127
$row8[0,5] = 0
$row8[0,6] = 2.169141230402332
$row8[0,7] = 0
$row8[0,8] = 2.903146184406467
$row8[0,9] = 0.1905340366506181
$row8[0,10] = 0
$row8[0,11] = 1.691143546854436
$ws.Range("C8:N8").Value = $row8

$row9 = New-Object 'object[,]' 1,12
$row9[0,0] = 0.1145635357669761
$row9[0,1] = 0.11696119741217
$row9[0,2] = 0.123941970925145
$row9[0,3] = 2.381942618816609
$row9[0,4] = 0.002536284982668457
$row9[0,5] = 0
$row9[0,6] = 2.219994697724019
$row9[0,7] = 0
$row9[0,8] = 3.360643847489598
$row9[0,9] = 0.2023815091552308
$row9[0,10] = 0
$row9[0,11] = 1.636153036063742
$ws.Range("C9:N9").Value = $row9

$row10 = New-Object 'object[,]' 1,12
$row10[0,0] = 0.1193278900634027
$row10[0,1] = 0.1146363081512227
$row10[0,2] = 0.1294507427361538
$row10[0,3] = 2.439236755446984
$row10[0,4] = 0.00252734623243504
$row10[0,5] = 0
$row10[0,6] = 2.264419693497402
$row10[0,7] = 0
$row10[0,8] = 3.704815571174947
$row10[0,9] = 0.2116900613961263
$row10[0,10] = 0
$row10[0,11] = 1.599092712011704
$ws.Range("C10:N10").Value = $row10

$row11 = New-Object 'object[,]' 1,12
$row11[0,0] = 0.1215458184955622
$row11[0,1] = 0.1136322946249528
$row11[0,2] = 0.132030118921449
$row11[0,3] = 2.467042611290054
$row11[0,4] = 0.002523466863089817
$row11[0,5] = 0
$row11[0,6] = 2.286191037646915
$row11[0,7] = 0
$row11[0,8] = 3.863190170111864
$row11[0,9] = 0.2160588261383936
$row11[0,10] = 0
$row11[0,11] = 1.582962238286479
$ws.Range("C11:N11").Value = $row11

$row12 = New-Object 'object[,]' 1,12
$row12[0,0] = 0.1223929767785421
$row12[0,1] = 0.1132598440845989
$row12[0,2] = 0.1330175054783354
$row12[0,3] = 2.477825125526749
$row12[0,4] = 0.002522024550404834
$row12[0,5] = 0
$row12[0,6] = 2.294662072709329
$row12[0,7] = 0
$row12[0,8] = 3.923426389383508
$row12[0,9] = 0.2177326837074816
$row12[0,10] = 0
$row12[0,11] = 1.57695917543467
$ws.Range("C12:N12").Value = $row12

$row13 = New-Object 'object[,]' 1,12
$row13[0,0] = 0.1222102024148484
$row13[0,1] = 0.1133397127248266
$row13[0,2] = 0.1328043799773582
$row13[0,3] = 2.475491624281574
$row13[0,4] = 0.002522333992243807
$row13[0,5] = 0
$row13[0,6] = 2.292827568066258
$row13[0,7] = 0
$row13[0,8] = 3.910441686693218
$row13[0,9] = 0.2173713180628738
$row13[0,10] = 0
$row13[0,11] = 1.578247356498677
$ws.Range("C13:N13").Value = $row13

$row14 = New-Object 'object[,]' 1,12
$row14[0,0] = 0.1216153688678361
$row14[0,1] = 0.1136014973802304
$row14[0,2] = 0.1321111382639728
$row14[0,3] = 2.467924609605149
$row14[0,4] = 0.002523347668433097
$row14[0,5] = 0
$row14[0,6] = 2.28688339949143
$row14[0,7] = 0
$row14[0,8] = 3.868140547543419
$row14[0,9] = 0.216196143443355
$row14[0,10] = 0
$row14[0,11] = 1.582466251780186
$ws.Range("C14:N14").Value = $row14

$row15 = New-Object 'object[,]' 1,12
$row15[0,0] = 0.121251963860729
$row15[0,1] = 0.1137628582522225
$row15[0,2] = 0.1316878947124138
$row15[0,3] = 2.46332262206036
$row15[0,4] = 0.002523972049310813
$row15[0,5] = 0
$row15[0,6] = 2.283272009252826
$row15[0,7] = 0
$row15[0,8] = 3.84226424607408
$row15[0,9] = 0.2154788605431435
$row15[0,10] = 0
$row15[0,11] = 1.585064160226574
$ws.Range("C15:N15").Value = $row15

$row16 = New-Object 'object[,]' 1,12
$row16[0,0] = 0.119183961451597
$row16[0,1] = 0.1147030046810791
$row16[0,2] = 0.1292836581128824
$row16[0,3] = 2.437454850962581
$row16[0,4] = 0.002527603506140874
$row16[0,5] = 0
$row16[0,6] = 2.263028495401514
$row16[0,7] = 0
$row16[0,8] = 3.694502081445933
$row16[0,9] = 0.2114072705406045
$row16[0,10] = 0
$row16[0,11] = 1.600161572675887
$ws.Range("C16:N16").Value = $row16

$row17 = New-Object 'object[,]' 1,12
$row17[0,0] = 0.1179282688474217
$row17[0,1] = 0.115293508916789
$row17[0,2] = 0.1278275924167218
$row17[0,3] = 2.422033847509368
$row17[0,4] = 0.002529879047886084
$row17[0,5] = 0
$row17[0,6] = 2.251011323607301
$row17[0,7] = 0
$row17[0,8] = 3.604320283248512
$row17[0,9] = 0.2089440061323415
$row17[0,10] = 0
$row17[0,11] = 1.609610234813362
$ws.Range("C17:N17").Value = $row17

$row18 = New-Object 'object[,]' 1,12
$row18[0,0] = 0.1172107902304731
$row18[0,1] = 0.1156381929659673
$row18[0,2] = 0.1269970067587423
$row18[0,3] = 2.413327952937721
$row18[0,4] = 0.002531205481617197
$row18[0,5] = 0
$row18[0,6] = 2.244246279898903
$row18[0,7] = 0
$row18[0,8] = 3.552620197068222
$row18[0,9] = 0.2075398266586461
$row18[0,10] = 0
$row18[0,11] = 1.615113393461654
$ws.Range("C18:N18").Value = $row18

$row19 = New-Object 'object[,]' 1,12
$row19[0,0] = 0.1169686822743898
$row19[0,1] = 0.1157557617231895
$row19[0,2] = 0.1267169683303848
$row19[0,3] = 2.410408354707386
$row19[0,4] = 0.002531657617707268
$row19[0,5] = 0
$row19[0,6] = 2.241980923064361
$row19[0,7] = 0
$row19[0,8] = 3.535144559660125
$row19[0,9] = 0.2070665581023263
$row19[0,10] = 0
$row19[0,11] = 1.616988426466856
$ws.Range("C19:N19").Value = $row19

$row20 = New-Object 'object[,]' 1,12
$row20[0,0] = 0.1180614465616401
$row20[0,1] = 0.1152301266470062
$row20[0,2] = 0.127981877923844
$row20[0,3] = 2.423658463092664
$row20[0,4] = 0.00252963499186802
$row20[0,5] = 0
$row20[0,6] = 2.252275352049551
$row20[0,7] = 0
$row20[0,8] = 3.613902663283341
$row20[0,9] = 0.2092049167875132
$row20[0,10] = 0
$row20[0,11] = 1.608597311970222
$ws.Range("C20:N20").Value = $row20

$row21 = New-Object 'object[,]' 1,12
$row21[0,0] = 0.1217898883935931
$row21[0,1] = 0.1135243942786204
$row21[0,2] = 0.132314470967259
$row21[0,3] = 2.470140337878206
$row21[0,4] = 0.002523049203218536
$row21[0,5] = 0
$row21[0,6] = 2.288623176752665
$row21[0,7] = 0
$row21[0,8] = 3.880558249757883
$row21[0,9] = 0.2165407898856415
$row21[0,10] = 0
$row21[0,11] = 1.581224201319809
$ws.Range("C21:N21").Value = $row21

$row22 = New-Object 'object[,]' 1,12
$row22[0,0] = 0.1242690574948568
$row22[0,1] = 0.1124547784317294
$row22[0,2] = 0.1352080840333798
$row22[0,3] = 2.501995018793565
$row22[0,4] = 0.002518900684123534
$row22[0,5] = 0
$row22[0,6] = 2.31370093979811
$row22[0,7] = 0
$row22[0,8] = 4.056368673492045
$row22[0,9] = 0.2214489681780805
$row22[0,10] = 0
$row22[0,11] = 1.563947678034959
$ws.Range("C22:N22").Value = $row22

$row23 = New-Object 'object[,]' 1,12
$row23[0,0] = 0.1229419965298035
$row23[0,1] = 0.1130215047186702
$row23[0,2] = 0.1336580081834242
$row23[0,3] = 2.48485768546675
$row23[0,4] = 0.002521100634674264
$row23[0,5] = 0
$row23[0,6] = 2.300194780635792
$row23[0,7] = 0
$row23[0,8] = 3.962393807042815
$row23[0,9] = 0.2188189064070656
$row23[0,10] = 0
$row23[0,11] = 1.573112208072745
$ws.Range("C23:N23").Value = $row23

$row24 = New-Object 'object[,]' 1,12
$row24[0,0] = 0.1180012231001371
$row24[0,1] = 0.1152587655919355
$row24[0,2] = 0.1279121051250343
$row24[0,3] = 2.422923476690499
$row24[0,4] = 0.002529745272861377
$row24[0,5] = 0
$row24[0,6] = 2.251703437077211
$row24[0,7] = 0
$row24[0,8] = 3.609570012775237
$row24[0,9] = 0.2090869217671099
$row24[0,10] = 0
$row24[0,11] = 1.60905503328175
$ws.Range("C24:N24").Value = $row24

$row25 = New-Object 'object[,]' 1,12
$row25[0,0] = 0.1128538671830057
$row25[0,1] = 0.1178646123491216
$row25[0,2] = 0.1219778820073891
$row25[0,3] = 2.362371994016996
$row25[0,4] = 0.002539743557945462
$row25[0,5] = 0
$row25[0,6] = 2.205006205943661
$row25[0,7] = 0
$row25[0,8] = 3.235482562891775
$row25[0,9] = 0.199071277086361
$row25[0,10] = 0
$row25[0,11] = 1.650443924827199
$ws.Range("C25:N25").Value = $row25
